$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = [double]"22.29000000000005"
$ws.Range("H2").Value = [double]"3.552713678800501e-16"
$ws.Range("K2").Value = [double]"43.01513842317185"
$ws.Range("L2").Value = "[32.88363982642909, 53.14663701991461]"
$ws.Range("M2").Value = [double]"2.886579864025407e-15"
$ws.Range("N2").Value = [double]"2.886579864025407e-15"
$ws.Range("O2").Value = [double]"1.276763380738195"
$ws.Range("P2").Value = "[1.0126054398958093, 1.54092132158058]"
$ws.Range("S2").Value = [double]"62.27009092323544"
$ws.Range("T2").Value = "[56.478162058375986, 68.06201978809489]"
$ws.Range("W2").Value = [double]"17.76060060060064"
$ws.Range("X2").Value = [double]"16.82348348348352"
$ws.Range("Y2").Value = [double]"18.69771771771775"

# Row 3
$ws.Range("E3").Value = [double]"22.53000000000008"
$ws.Range("H3").Value = [double]"3.552713678800501e-16"
$ws.Range("K3").Value = [double]"38.81316328328045"
$ws.Range("L3").Value = "[30.919812235304107, 46.706514331256784]"
$ws.Range("M3").Value = [double]"0"
$ws.Range("N3").Value = [double]"0"
$ws.Range("O3").Value = [double]"1.490605523324887"
$ws.Range("P3").Value = "[1.2641844311742707, 1.7170266154755032]"
$ws.Range("Q3").Value = [double]"0"
$ws.Range("R3").Value = [double]"0"
$ws.Range("S3").Value = [double]"58.47733435733272"
$ws.Range("T3").Value = "[53.43085402857188, 63.523814686093566]"
$ws.Range("W3").Value = [double]"17.18504504504511"
$ws.Range("X3").Value = [double]"16.37315315315321"
$ws.Range("Y3").Value = [double]"17.996936936937"
